$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at H (shifts old H:R -> I:S)
$ws.Columns("H").Insert()

# Update header row
$ws.Range("G1").Value2 = "norm_k"
$ws.Range("H1").Value2 = "k_Truss"

# Row 2 (Graph1)
$ws.Range("C2").Value2 = 3
$ws.Range("H2").Value2 = 3
$ws.Range("I2").Value2 = -0.03777544596012592
$ws.Range("J2").Value2 = 0.1934321414709257
$ws.Range("K2").Value2 = 0.1857457212598542
$ws.Range("L2").Value2 = -0.03777544596012591
$ws.Range("M2").Value2 = 0.4166666666666667
$ws.Range("N2").Value2 = 1.887918502671133
$ws.Range("O2").Value2 = 0.6500224216483541
$ws.Range("P2").Value2 = 0.2261904761904762
$ws.Range("Q2").Value2 = 0.4040610178208843
$ws.Range("R2").Value2 = "Subgraph 1 (Nodes): [1, 0, 16, 19] - Density: 1`nSubgraph 2 (Nodes): [4, 24, 15, 13] - Density: 1`nSubgraph 3 (Nodes): [10, 6, 5, 25, 11, 22] - Density: 0.799943`n----------------------------------------------------"
$ws.Range("S2").Value2 = "--- DEBUG: Seeds Loaded ---`nTotal seeds = 2`nSeed 1: { 5 6 10 11 22 25 } | Triangles: 16 | Density: 0.799943`nSeed 2: { 4 13 15 24 } | Triangles: 4 | Density: 1`nSubgraph 1:z { 0 1 2 4 5 6 10 11 13 15 16 17 19 20 22 23 24 25 26 27 } N: 20 Triangles: 25 Density: 0.0219297"

# Row 3 (Graph2)
$ws.Range("C3").Value2 = 3
$ws.Range("D3").Value2 = 0.45
$ws.Range("H3").Value2 = 3
$ws.Range("I3").Value2 = 0
$ws.Range("J3").Value2 = 0
$ws.Range("K3").Value2 = 0
$ws.Range("L3").Value2 = 0
$ws.Range("M3").Value2 = 0.35
$ws.Range("N3").Value2 = 1.93703267660925
$ws.Range("O3").Value2 = 0
$ws.Range("P3").Value2 = 0.1296296296296296
$ws.Range("Q3").Value2 = 0.5025575614435649
$ws.Range("R3").Value2 = "Subgraph 1 (Nodes): [46, 45, 19, 43, 14, 35, 9, 41, 28, 25, 4, 42, 27, 5] - Density: 0.513739`nSubgraph 2 (Nodes): [36, 8, 33, 40, 11, 48] - Density: 0.549961`nSubgraph 3 (Nodes): [39, 26, 23, 15, 3, 1, 7, 47, 6] - Density: 0.607122`n----------------------------------------------------"
$ws.Range("S3").Value2 = "--- DEBUG: Seeds Loaded ---`nTotal seeds = 3`nSeed 1: { 8 11 33 36 40 48 } | Triangles: 11 | Density: 0.549961`nSeed 2: { 4 5 9 14 19 25 27 28 35 41 42 43 45 46 } | Triangles: 187 | Density: 0.513739`nSeed 3: { 1 3 6 7 15 23 26 39 47 } | Triangles: 51 | Density: 0.607122`nSubgraph 1:z { 1 2 3 4 5 6 7 8 9 10 11 14 15 19 20 21 22 23 25 26 27 28 29 30 32 33 34 35 36 38 39 40 41 42 43 44 45 46 47 48 } N: 40 Triangles: 249 Density: 0.0251435"

# Row 4 (Graph1, duplicate of row 2)
$ws.Range("C4").Value2 = 3
$ws.Range("H4").Value2 = 3
$ws.Range("I4").Value2 = -0.03777544596012592
$ws.Range("J4").Value2 = 0.1934321414709257
$ws.Range("K4").Value2 = 0.1857457212598542
$ws.Range("L4").Value2 = -0.03777544596012591
$ws.Range("M4").Value2 = 0.4166666666666667
$ws.Range("N4").Value2 = 1.887918502671133
$ws.Range("O4").Value2 = 0.6500224216483541
$ws.Range("P4").Value2 = 0.2261904761904762
$ws.Range("Q4").Value2 = 0.4040610178208843
$ws.Range("R4").Value2 = "Subgraph 1 (Nodes): [1, 0, 16, 19] - Density: 1`nSubgraph 2 (Nodes): [4, 24, 15, 13] - Density: 1`nSubgraph 3 (Nodes): [10, 6, 5, 25, 11, 22] - Density: 0.799943`n----------------------------------------------------"
$ws.Range("S4").Value2 = "--- DEBUG: Seeds Loaded ---`nTotal seeds = 2`nSeed 1: { 5 6 10 11 22 25 } | Triangles: 16 | Density: 0.799943`nSeed 2: { 4 13 15 24 } | Triangles: 4 | Density: 1`nSubgraph 1:z { 0 1 2 4 5 6 10 11 13 15 16 17 19 20 22 23 24 25 26 27 } N: 20 Triangles: 25 Density: 0.0219297"

# Row 5 (Graph2, duplicate of row 3)
$ws.Range("C5").Value2 = 3
$ws.Range("D5").Value2 = 0.45
$ws.Range("H5").Value2 = 3
$ws.Range("I5").Value2 = 0
$ws.Range("J5").Value2 = 0
$ws.Range("K5").Value2 = 0
$ws.Range("L5").Value2 = 0
$ws.Range("M5").Value2 = 0.35
$ws.Range("N5").Value2 = 1.93703267660925
$ws.Range("O5").Value2 = 0
$ws.Range("P5").Value2 = 0.1296296296296296
$ws.Range("Q5").Value2 = 0.5025575614435649
$ws.Range("R5").Value2 = "Subgraph 1 (Nodes): [46, 45, 19, 43, 14, 35, 9, 41, 28, 25, 4, 42, 27, 5] - Density: 0.513739`nSubgraph 2 (Nodes): [36, 8, 33, 40, 11, 48] - Density: 0.549961`nSubgraph 3 (Nodes): [39, 26, 23, 15, 3, 1, 7, 47, 6] - Density: 0.607122`n----------------------------------------------------"
$ws.Range("S5").Value2 = "--- DEBUG: Seeds Loaded ---`nTotal seeds = 3`nSeed 1: { 8 11 33 36 40 48 } | Triangles: 11 | Density: 0.549961`nSeed 2: { 4 5 9 14 19 25 27 28 35 41 42 43 45 46 } | Triangles: 187 | Density: 0.513739`nSeed 3: { 1 3 6 7 15 23 26 39 47 } | Triangles: 51 | Density: 0.607122`nSubgraph 1:z { 1 2 3 4 5 6 7 8 9 10 11 14 15 19 20 21 22 23 25 26 27 28 29 30 32 33 34 35 36 38 39 40 41 42 43 44 45 46 47 48 } N: 40 Triangles: 249 Density: 0.0251435"
